# Add a new "Croatia" market test-data sheet, based on the existing
# "Slovakia" sheet (same column widths / layout), trimmed down to the
# 10-row "Turkey"-style layout and filled in with the Croatia values.

$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("Slovakia")

# Select the whole source sheet (mirrors the copy/duplicate workflow the
# author used, and is what leaves the "select all" selection behind on
# the Slovakia sheet afterwards).
$source.Activate() | Out-Null
$source.Cells.Select() | Out-Null

# Duplicate it to the end of the workbook.
$source.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Croatia"

# The Slovakia sheet has two extra product rows (XLM800-STI / XLM800-Zetfas)
# that Croatia doesn't need - drop them so the sheet matches the 10-row
# layout used by the other markets.
$newSheet.Rows("9:10").Delete()

# Fill in the market-specific values.
$newSheet.Range("B2").Value = "Croatia Market"
$newSheet.Range("B4").Value = "NGC-3139/T2473"

# Leave the same kind of cell selected/active as the source sheet had.
$newSheet.Range("G16").Select() | Out-Null
